# Leave Card update 5/18/2023 4:35 PM
#
# The document contains two copies of the same "Certificate of Leave
# Credits" (cached MERGEFIELD results). Update each field result in place.
#
# Most fields are simple document-wide Find/Replace calls. The "day of
# month" (31) and its superscript ordinal suffix (st) sit in two adjacent
# runs with no separating character ("31" immediately followed by "st"),
# so a whole-document / whole-word Find can't isolate "st" from the
# unrelated word "request" elsewhere in the same sentence. Those two are
# therefore replaced with the Find scoped to just their own paragraph.

$d = $word.ActiveDocument

function Replace-All($find, $replace, $wholeWord) {
    $d.Content.Find.Execute($find, $true, $wholeWord, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

function Replace-InParagraph($paraIndex, $find, $replace) {
    $rng = $d.Paragraphs.Item($paraIndex).Range
    $rng.Find.Execute($find, $true, $false, $false, $false, $false, `
                       $true, 1, $false, $replace, 2) | Out-Null
}

# Salutation (caps form, e.g. "MR" -> "MS")
Replace-All "MR" "MS" $true

# Full name
Replace-All "FRANZ R. VIDA" "MA. CARMELA D. ARELLANO" $false

# Position
Replace-All "Casual Employee" "Casual Medical Technologist" $false

# Office
Replace-All "Vice Mayor's Office" "Ospital Ng Tagaytay" $false

# Salutation (title case form, e.g. "Mr" -> "Ms")
Replace-All "Mr" "Ms" $true

# Last name
Replace-All "Vida" "Arellano" $true

# "Issued this 31st day of March 2023" paragraphs (two copies on the page)
foreach ($idx in 14, 36) {
    Replace-InParagraph $idx "31" "18"
    Replace-InParagraph $idx "st" "th"
    Replace-InParagraph $idx "March" "May"
}

# Last day of service (must run after the paragraph-scoped "31" -> "18"
# above, otherwise the "31" introduced here would get clobbered too)
Replace-All "June 30, 2022" "January 31, 2022" $false

# Leave credit figures (note the exact leading spaces from the field formatting)
Replace-All "  47.500" "  41.250" $false
Replace-All "  65.500" "  61.250" $false
Replace-All " 113.000" " 102.500" $false
